# Add a new "Greece" market sheet, cloned from the existing "Croatia" sheet
# (same layout/styles), populated with the Greece-specific values, and
# placed immediately after Croatia as the new active tab.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Clone Croatia (keeps styles, merged cells, column widths, etc.) and place
# the copy right after it.
$croatia.Copy($null, $croatia)
$greece = $wb.Worksheets.Item($croatia.Index + 1)
$greece.Name = "Greece"

# Fill in the Greece-specific values (NGC ticket reference first, then the
# market name, matching the original authoring order).
$greece.Range("B4").Value = "NGC-4119/T3187/T3189/T3185"
$greece.Range("B2").Value = "Greece Market"

# Restore Croatia's selection to a "whole sheet" state (it is no longer the
# active tab) and make Greece the active tab with its own selection.
$croatia.Activate()
$croatia.Cells.Select() | Out-Null

$greece.Activate()
$greece.Range("C14").Select() | Out-Null
